$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("J8").Value  = 0.1485032540344368
$ws.Range("I9").Value  = 0.1734537503564907
$ws.Range("H10").Value = 0.2422520263583712
$ws.Range("G11").Value = 0.2534537503564908
$ws.Range("F12").Value = 0.381103329907261
$ws.Range("E13").Value = 0.04235042473292953
$ws.Range("D14").Value = 0.07961008106920435
$ws.Range("C15").Value = 0.02893023050567838
$ws.Range("B16").Value = 0.02940328597706714
